# Apply weekly cryptocurrency market-data refresh to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.416.23'
$ws.Range('E2').Value = '  +4.35%  '
$ws.Range('D3').Value = '3.489.58'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.73'
$ws.Range('E5').Value = '  +2.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.76'
$ws.Range('E6').Value = '  +7.48%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.70'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +4.67%  '
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('D12').Value = '4.087.80'
$ws.Range('E12').Value = '  +3.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.72'
$ws.Range('E13').Value = '  +6.44%  '
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '3.489.79'
$ws.Range('E15').Value = '  +3.63%  '
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('D17').Value = '63.410.01'
$ws.Range('E17').Value = '  +4.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.33'
$ws.Range('E18').Value = '  +4.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.37'
$ws.Range('E19').Value = '  +6.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.39'
$ws.Range('E20').Value = '  +5.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '393.22'
$ws.Range('E21').Value = '  +2.76%  '
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '75.31'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +8.83%  '
$ws.Range('D26').Value = '3.629.00'
$ws.Range('E26').Value = '  +3.65%  '
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.86'
$ws.Range('E28').Value = '  +10.59%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +5.14%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.15'
$ws.Range('E31').Value = '  +2.51%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.44'
$ws.Range('E32').Value = '  +7.41%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.84'
$ws.Range('E34').Value = '  +3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '32.65'
$ws.Range('E35').Value = '  +27.23%  '
$ws.Range('E36').Value = '  +8.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.16'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '171.99'
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('E39').Value = '  +9.79%  '
$ws.Range('D40').Value = '3.526.73'
$ws.Range('E40').Value = '  +3.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0766'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.806'
$ws.Range('E42').Value = '  +4.81%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.74'
$ws.Range('E43').Value = '  +7.65%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.50'
$ws.Range('E44').Value = '  +4.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.50'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.22'
$ws.Range('E46').Value = '  +10.35%  '
$ws.Range('D47').Value = '2.623.83'
$ws.Range('E47').Value = '  +7.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.82'
$ws.Range('E48').Value = '  +8.00%  '
$ws.Range('E49').Value = '  +17.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.76'
$ws.Range('E50').Value = '  +2.22%  '
$ws.Range('E51').Value = '  +5.26%  '

Write-Output "Applied 88 cell updates to cryptos sheet"
